# Timesheet update for 2019-04-14.xlsx ("timesheet.hourly" sheet)
# Fills in the "Reg. hours" row (row 12) for several days of the
# second week, matching the author's "Updated App and Firmware" commit:
#   E12 (MON) = 3.5
#   G12 (WED) = 3.75
#   H12 (THU) = 1.5
#   L12 (MON, wk2) = 4.5
#   N12 (WED, wk2) = 4.25
# The weekly-total column (R12) and the summary row (row 22) are driven
# by existing formulas, so they recompute automatically on recalc.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("timesheet.hourly")

$ws.Range("E12").Value = 3.5
$ws.Range("G12").Value = 3.75
$ws.Range("H12").Value = 1.5
$ws.Range("L12").Value = 4.5
$ws.Range("N12").Value = 4.25

# Recalculate so dependent totals (R12, row 22) reflect the new inputs.
$excel.Calculate()

# Move the active selection to H13, matching the saved cursor position.
$ws.Range("H13").Select() | Out-Null
